$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. '230.40', '17.00').
# Force text format before assignment so Excel keeps them as literal text
# (matching the original inline-string cells) instead of coercing to numbers,
# then restore the Normal style so no stray 's' attribute sticks to the cell.

$dCell = $ws.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "38.189.11"
$dCell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +2.88%  "

$dCell = $ws.Cells.Item(3, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.060.74"
$dCell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +2.37%  "

$ws.Cells.Item(4, 5).Value = "  -0.52%  "

$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "230.40"
$dCell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +1.85%  "

$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.615"
$dCell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.41%  "

$dCell = $ws.Cells.Item(7, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "59.73"
$dCell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +8.12%  "

$dCell = $ws.Cells.Item(9, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.387"
$dCell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +3.19%  "

$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0812"
$dCell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +4.06%  "

$ws.Cells.Item(11, 5).Value = "  +2.35%  "

$ws.Cells.Item(12, 5).Value = "  +5.07%  "

$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.366.84"
$dCell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +2.43%  "

$dCell = $ws.Cells.Item(14, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "21.33"
$dCell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +7.53%  "

$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.756"
$dCell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +2.46%  "

$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "5.31"
$dCell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.89%  "

$dCell = $ws.Cells.Item(17, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.063.52"
$dCell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +2.45%  "

$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "38.110.78"
$dCell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +2.87%  "

$dCell = $ws.Cells.Item(19, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.27"
$dCell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.84%  "

$dCell = $ws.Cells.Item(20, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "69.96"
$dCell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.31%  "

$dCell = $ws.Cells.Item(21, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0836"
$dCell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +2.78%  "

$dCell = $ws.Cells.Item(22, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "225.41"
$dCell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.72%  "

$ws.Cells.Item(23, 5).Value = "  +0.02%  "

$ws.Cells.Item(24, 5).Value = "  +0.13%  "

$ws.Cells.Item(25, 5).Value = "  +4.26%  "

$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "9.33"
$dCell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +4.23%  "

$dCell = $ws.Cells.Item(27, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "166.24"
$dCell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +1.04%  "

$ws.Cells.Item(28, 5).Value = "  +7.15%  "

$dCell = $ws.Cells.Item(29, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "19.06"
$dCell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.28%  "

$ws.Cells.Item(30, 5).Value = "  +2.89%  "

$dCell = $ws.Cells.Item(31, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.119"
$dCell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +2.37%  "

$dCell = $ws.Cells.Item(32, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "4.57"
$dCell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +2.53%  "

$ws.Cells.Item(33, 5).Value = "  +2.94%  "

$ws.Cells.Item(34, 5).Value = "  +10.77%  "

$ws.Cells.Item(35, 5).Value = "  +1.36%  "

$ws.Cells.Item(36, 5).Value = "  +1.24%  "

$dCell = $ws.Cells.Item(37, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.12"
$dCell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +14.98%  "

$ws.Cells.Item(38, 5).Value = "  +5.59%  "

$ws.Cells.Item(39, 5).Value = "  -0.09%  "

$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.533.04"
$dCell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +5.29%  "

$dCell = $ws.Cells.Item(41, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "98.48"
$dCell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +3.92%  "

$ws.Cells.Item(42, 5).Value = "  +2.64%  "

$dCell = $ws.Cells.Item(43, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "17.00"
$dCell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +5.76%  "

$ws.Cells.Item(44, 5).Value = "  +4.41%  "

$dCell = $ws.Cells.Item(45, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0923"
$dCell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.98%  "

$ws.Cells.Item(46, 5).Value = "  +1.13%  "

$dCell = $ws.Cells.Item(47, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "4.12"
$dCell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +9.59%  "

$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.02"
$dCell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +2.60%  "

$dCell = $ws.Cells.Item(49, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.97"
$dCell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +2.85%  "

$dCell = $ws.Cells.Item(50, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "7.13"
$dCell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.07%  "

$dCell = $ws.Cells.Item(51, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.254.18"
$dCell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +2.53%  "
